$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("E2").Value = 11.7057952
$ws.Range("E3").Value = 2.1799264
$ws.Range("E4").Value = 0.885464
$ws.Range("E5").Value = 1.419984
$ws.Range("E6").Value = 0.9508511999999999
$ws.Range("E7").Value = 0.2785056
$ws.Range("E8").Value = 0.08363520000000001
$ws.Range("C9").Value = 843
$ws.Range("E9").Value = 0.8668096000000001
$ws.Range("C10").Value = 598
$ws.Range("E10").Value = 0.5152896
$ws.Range("C11").Value = 1812
$ws.Range("E11").Value = 0.758016
$ws.Range("C12").Value = 411
$ws.Range("E12").Value = 1.2240096
$ws.Range("C13").Value = 3241
$ws.Range("E13").Value = 12.58827199999999
$ws.Range("E14").Value = 4.864020800000001
$ws.Range("C15").Value = 290
$ws.Range("E15").Value = 0.8643535999999999
$ws.Range("C16").Value = 419
$ws.Range("E16").Value = 0.8735807999999999
$ws.Range("C17").Value = 603
$ws.Range("E17").Value = 1.387152
$ws.Range("C18").Value = 175
$ws.Range("E18").Value = 0.436464
$ws.Range("C19").Value = 57
$ws.Range("E19").Value = 0.1157376
$ws.Range("E20").Value = 62.82544
$ws.Range("E21").Value = 66.695296
$ws.Range("E22").Value = 82.556928
$ws.Range("E23").Value = 254.024992

$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("C2").Value = 602
$ws.Range("D2").Value = 141
$ws.Range("F2").Value = 294
$ws.Range("G2").Value = 82
$ws.Range("H2").Value = 24
$ws.Range("I2").Value = 420
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 904
$ws.Range("L2").Value = 208
$ws.Range("C3").Value = 606
$ws.Range("D3").Value = 147
$ws.Range("F3").Value = 294
$ws.Range("G3").Value = 86
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 423
$ws.Range("J3").Value = 298
$ws.Range("K3").Value = 908
$ws.Range("L3").Value = 203
$ws.Range("C4").Value = 602
$ws.Range("D4").Value = 142
$ws.Range("F4").Value = 301
$ws.Range("G4").Value = 87
$ws.Range("H4").Value = 30
$ws.Range("J4").Value = 297
$ws.Range("K4").Value = 903
$ws.Range("L4").Value = 208
$ws.Range("C5").Value = 594
$ws.Range("D5").Value = 148
$ws.Range("F5").Value = 302
$ws.Range("G5").Value = 88
$ws.Range("H5").Value = 27
$ws.Range("I5").Value = 422
$ws.Range("J5").Value = 297
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 211
$ws.Range("C6").Value = 597
$ws.Range("D6").Value = 155
$ws.Range("F6").Value = 291
$ws.Range("G6").Value = 90
$ws.Range("H6").Value = 25
$ws.Range("I6").Value = 421
$ws.Range("J6").Value = 299
$ws.Range("K6").Value = 897
$ws.Range("L6").Value = 208

$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("C2").Value = 75.25
$ws.Range("D2").Value = 17.625
$ws.Range("F2").Value = 36.75
$ws.Range("G2").Value = 10.25
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 52.5
$ws.Range("J2").Value = 37.5
$ws.Range("K2").Value = 113
$ws.Range("L2").Value = 26
$ws.Range("C3").Value = 113.625
$ws.Range("D3").Value = 27.5625
$ws.Range("F3").Value = 55.125
$ws.Range("G3").Value = 16.125
$ws.Range("H3").Value = 4.875
$ws.Range("I3").Value = 79.3125
$ws.Range("J3").Value = 55.875
$ws.Range("K3").Value = 170.25
$ws.Range("L3").Value = 38.0625
$ws.Range("C4").Value = 131.6875
$ws.Range("D4").Value = 31.0625
$ws.Range("F4").Value = 65.84375
$ws.Range("G4").Value = 19.03125
$ws.Range("H4").Value = 6.5625
$ws.Range("J4").Value = 64.96875
$ws.Range("K4").Value = 197.53125
$ws.Range("L4").Value = 45.5
$ws.Range("C5").Value = 139.21875
$ws.Range("D5").Value = 34.6875
$ws.Range("F5").Value = 70.78125
$ws.Range("G5").Value = 20.625
$ws.Range("H5").Value = 6.328125
$ws.Range("I5").Value = 98.90625
$ws.Range("J5").Value = 69.609375
$ws.Range("K5").Value = 210.9375
$ws.Range("L5").Value = 49.453125
$ws.Range("C6").Value = 144.5859375
$ws.Range("D6").Value = 37.5390625
$ws.Range("F6").Value = 70.4765625
$ws.Range("G6").Value = 21.796875
$ws.Range("H6").Value = 6.0546875
$ws.Range("I6").Value = 101.9609375
$ws.Range("J6").Value = 72.4140625
$ws.Range("K6").Value = 217.2421875
$ws.Range("L6").Value = 50.375

$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B3").Value = 522550
$ws.Range("B4").Value = 522550
$ws.Range("B5").Value = 4926899.999999999
